$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''51.736.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '''2.799.20'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.76%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''353.16'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').Value = '''111.78'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.54%  '
$ws.Range('D7').Value = '''0.557'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.55%  '
$ws.Range('D9').Value = '''0.625'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.20%  '
$ws.Range('D10').Value = '''40.21'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.05%  '
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('D12').Value = '''0.0837'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').Value = '''19.93'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('D14').Value = '''7.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.73%  '
$ws.Range('D15').Value = '''3.231.43'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.61%  '
$ws.Range('D16').Value = '''2.831.49'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.17%  '
$ws.Range('D17').Value = '''0.945'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.02%  '
$ws.Range('D18').Value = '''51.690.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('D19').Value = '''7.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('E20').Value = '  +5.77%  '
$ws.Range('E21').Value = '  +4.55%  '
$ws.Range('D22').Value = '''0.0₃0972'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('D23').Value = '''70.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.40%  '
$ws.Range('D24').Value = '''267.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.04%  '
$ws.Range('D25').Value = '''2.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '''26.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('D28').Value = '''0.160'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').Value = '''38.95'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +13.21%  '
$ws.Range('D30').Value = '''10.37'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.21%  '
$ws.Range('D31').Value = '''2.27'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('D32').Value = '''52.64'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.92%  '
$ws.Range('D33').Value = '''6.12'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.00%  '
$ws.Range('D34').Value = '''0.0453'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.85%  '
$ws.Range('D35').Value = '''0.0888'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.22%  '
$ws.Range('D36').Value = '''5.55'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.40%  '
$ws.Range('D37').Value = '''1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').Value = '''18.86'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.61%  '
$ws.Range('D39').Value = '''2.01'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.67%  '
$ws.Range('D40').Value = '''3.16'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('E41').Value = '  +1.96%  '
$ws.Range('D42').Value = '''2.51'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.68%  '
$ws.Range('E43').Value = '  +1.05%  '
$ws.Range('D44').Value = '''120.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.70%  '
$ws.Range('D45').Value = '''21.79'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '''2.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.64%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '''3.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.17%  '
$ws.Range('D48').Value = '''2.102.55'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('D49').Value = '''0.955'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.85%  '
$ws.Range('D50').Value = '''5.48'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').Value = '''1.37'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.00%  '
